$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a
# new ListBullet paragraph right after it, listing the two instructors
# on separate lines (joined by a manual line break, matching the style
# used elsewhere in the document for multi-line bullet paragraphs).

$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Docente(s)*") {
        $target = $para
        break
    }
}

$target.Range.InsertParagraphAfter()
$newIndex = $target.Index + 1

$firstPara = $d.Paragraphs($newIndex)
$firstPara.Style = "ListBullet"
$firstPara.Range.Text = "5111420 - Talita Martins Lacerda"

$firstPara = $d.Paragraphs($newIndex)
$firstPara.Range.InsertParagraphAfter()

$secondPara = $d.Paragraphs($newIndex + 1)
$secondPara.Style = "ListBullet"
$secondPara.Range.Text = "8853480 - Tatiane da Franca Silva"

# Merge the two paragraphs into one, turning the paragraph mark that
# separates them into a manual line break (<w:br/>) instead, so the two
# names live in a single ListBullet paragraph as two runs.
$firstPara = $d.Paragraphs($newIndex)
$markPos = $firstPara.Range.End - 1
$insertPoint = $d.Range($markPos, $markPos)
$insertPoint.InsertAfter([char]11)

$markRange = $d.Range($markPos + 1, $markPos + 2)
$markRange.Delete()

Write-Host "Inserted bullet paragraph with 2 instructors after 'Docente(s) Responsável(eis)'"
